# Apply the WR_89719272_WeekEnding_081025.xlsx edits:
#  - refresh the "Report Generated On" timestamp
#  - recompute summary totals (Total Billed Amount / Total Line Items)
#  - clear the Scope ID # value
#  - collapse the two detail line items (Point 02 PLA-CUT + Point 11 POL-40-2)
#    down to the single remaining line item (Point 11 POL-40-2), which pushes
#    the TOTAL row up one row and removes the now-unused gray banding row
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- header / summary fields -------------------------------------------------
$ws.Range("D5").Value = "Report Generated On: 08/26/2025 10:01 AM"
$ws.Range("C8").Value = 478.55
$ws.Range("C9").Value = 1
$ws.Range("G10").Value = ""

# --- overwrite row 16 with the surviving line item (was row 17) -------------
$ws.Range("A16").Value = "Point 11"
$ws.Range("B16").Value = "POL-40-2"
$ws.Range("C16").Value = "Inst"
$ws.Range("D16").Value = "Pole,40ft,Class 2"
$ws.Range("E16").Value = "EA"
$ws.Range("F16").Value = 1
$ws.Range("H16").Value = 478.55

# --- drop the old row 17 (its data now lives in row 16); TOTAL row shifts up
# to row 17, mergeCells/dimension are fixed up automatically by the delete.
$ws.Rows(17).Delete()

# --- TOTAL row (now row 17) ---------------------------------------------------
$ws.Range("H17").Value = 478.55
